# Fruta / hortaliza, semanal
# Rotate the weekly price-record data (columns D, L-T) across rows
# 2,3,4,5,6,7,8,9,10,12 of the sheet. Row 11 is left untouched, and the
# descriptive columns A,B,C,E-K are left untouched on every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters -> column numbers used below:
#   D=4, L=12, M=13, N=14, O=15, P=16, Q=17, R=18, S=19, T=20

function Get-RowData($r) {
    [ordered]@{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

function Set-RowData($r, $data) {
    $ws.Cells.Item($r, 4).Value2  = $data.D
    $ws.Cells.Item($r, 12).Value2 = $data.L
    $ws.Cells.Item($r, 13).Value2 = $data.M
    $ws.Cells.Item($r, 14).Value2 = $data.N
    $ws.Cells.Item($r, 15).Value2 = $data.O
    $ws.Cells.Item($r, 16).Value2 = $data.P
    $ws.Cells.Item($r, 17).Value2 = $data.Q
    $ws.Cells.Item($r, 18).Value2 = $data.R
    $ws.Cells.Item($r, 19).Value2 = $data.S
    $ws.Cells.Item($r, 20).Value2 = $data.T
}

# Snapshot the current ("before") data for every affected row first, since
# the rotation below would otherwise clobber values before they are read.
$snapshot = @{}
foreach ($r in 2,3,4,5,6,7,8,9,10,12) {
    $snapshot[$r] = Get-RowData $r
}

# Target row <- source row (old data that should now live in target row).
$rotation = @{
    2  = 4
    3  = 6
    4  = 8
    5  = 3
    6  = 9
    7  = 12
    8  = 7
    9  = 2
    10 = 5
    12 = 10
}

foreach ($target in $rotation.Keys) {
    $source = $rotation[$target]
    Set-RowData $target $snapshot[$source]
}
